$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Insert a new row above the old row 9 ("pastTimeHorizon") to host the new
# "round_for_capacity_market" parameter. Excel will auto-shift every row
# below (and rewrite dependent formula references) down by one.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "round_for_capacity_market"
$ws.Range("B9").Value = $false
$ws.Range("C9").Value = "don’t change this"

# investment_initialization_years (now row 19) changes from 0 to 3
$ws.Range("B19").Value = 3

# The conditional formatting range that used to cover B49:B54 needs to keep
# tracking the "ok" check rows, which are now B50:B55.
$cfRange = $ws.Range("B50:B55")
$cf = $cfRange.FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("B50:B55"))

# Match the author's final selection position.
$ws.Range("C9").Select()
